$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 909.25
$ws.Range("I19").Value = 824
$ws.Range("K19").Value = 824
$ws.Range("M19").Value = -649
$ws.Range("H28").Value = 361.4
$ws.Range("I28").Value = 278.08334
$ws.Range("K28").Value = 278.08334
$ws.Range("M28").Value = 206.91666
$ws.Range("H92").Value = 1309.7
$ws.Range("I92").Value = 2105.75
$ws.Range("J92").Value = 779
$ws.Range("K92").Value = 2105.75
$ws.Range("L92").Value = 779
$ws.Range("M92").Value = -857.75
$ws.Range("N92").Value = -3275
$ws.Range("H100").Value = 1165.5714
$ws.Range("I100").Value = 1084.8334
$ws.Range("J100").Value = 1650
$ws.Range("K100").Value = 1084.8334
$ws.Range("L100").Value = 1650
$ws.Range("M100").Value = -543.8334
$ws.Range("N100").Value = -2732
$ws.Range("H131").Value = 533102.5
$ws.Range("I131").Value = 674170
$ws.Range("K131").Value = 2022510
$ws.Range("M131").Value = -2017470
$ws.Range("H132").Value = 3798
$ws.Range("I132").Value = 4389.514
$ws.Range("K132").Value = 13168.542
$ws.Range("M132").Value = -10638.542
$ws.Range("H137").Value = 2001914.1
$ws.Range("I137").Value = 3127092.2
$ws.Range("K137").Value = 9381276.600000001
$ws.Range("M137").Value = -9378726.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1883620.5
$ws.Range("I32").Value = 905899.5600000001
$ws.Range("K32").Value = 905899.5600000001
$ws.Range("M32").Value = -905612.5600000001
$ws.Range("H63").Value = 1249.5
$ws.Range("I63").Value = 1249.5
$ws.Range("K63").Value = 1249.5
$ws.Range("M63").Value = -563.5
$ws.Range("H66").Value = 1249.5
$ws.Range("I66").Value = 1249.5
$ws.Range("K66").Value = 6247.5
$ws.Range("M66").Value = -2815.5
$ws.Range("H74").Value = 1411.12
$ws.Range("I74").Value = 1066.4445
$ws.Range("K74").Value = 1066.4445
$ws.Range("M74").Value = -192.4445000000001
$ws.Range("H77").Value = 1411.12
$ws.Range("I77").Value = 1066.4445
$ws.Range("K77").Value = 5332.2225
$ws.Range("M77").Value = -964.2224999999999
$ws.Range("H122").Value = 3899.4
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 2997
$ws.Range("N122").Value = -7897

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 13001684
$ws.Range("I105").Value = 1251523.8
$ws.Range("K105").Value = 1251523.8
$ws.Range("M105").Value = -1249776.8
$ws.Range("H134").Value = 2588.4443
$ws.Range("I134").Value = 2185.2144
$ws.Range("K134").Value = 6555.6432
$ws.Range("M134").Value = -4020.6432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 160.66667
$ws.Range("I7").Value = 110.545456
$ws.Range("K7").Value = 110.545456
$ws.Range("M7").Value = 2.454543999999999
$ws.Range("H15").Value = 3000
$ws.Range("J15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("N15").Value = -3340
$ws.Range("H31").Value = 7356974.5
$ws.Range("I31").Value = 3721.818
$ws.Range("K31").Value = 3721.818
$ws.Range("M31").Value = -3426.818
$ws.Range("H34").Value = 7356974.5
$ws.Range("I34").Value = 3721.818
$ws.Range("K34").Value = 3721.818
$ws.Range("M34").Value = -3519.818
$ws.Range("H50").Value = 59666
$ws.Range("J50").Value = 59666
$ws.Range("L50").Value = 59666
$ws.Range("N50").Value = -60916
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 1626
$ws.Range("I58").Value = 1145.8667
$ws.Range("K58").Value = 1145.8667
$ws.Range("M58").Value = -942.8667
$ws.Range("H59").Value = 75500
$ws.Range("J59").Value = 75500
$ws.Range("L59").Value = 75500
$ws.Range("N59").Value = -77790
$ws.Range("H60").Value = 19233.111
$ws.Range("J60").Value = 19887.25
$ws.Range("L60").Value = 19887.25
$ws.Range("N60").Value = -20909.25
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 85000
$ws.Range("J68").Value = 85000
$ws.Range("L68").Value = 85000
$ws.Range("N68").Value = -86498
$ws.Range("H71").Value = 85000
$ws.Range("J71").Value = 85000
$ws.Range("L71").Value = 255000
$ws.Range("N71").Value = -262488
$ws.Range("H99").Value = 5970
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 5970
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 5970
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -8966
$ws.Range("H126").Value = 5970
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5970
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 17910
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -22850
$ws.Range("H132").Value = 4115.9062
$ws.Range("I132").Value = 3400.36
$ws.Range("K132").Value = 10201.08
$ws.Range("M132").Value = -7671.08
$ws.Range("H134").Value = 3815.647
$ws.Range("I134").Value = 3890.862
$ws.Range("J134").Value = 3379.4
$ws.Range("K134").Value = 11672.586
$ws.Range("L134").Value = 10138.2
$ws.Range("M134").Value = -9137.585999999999
$ws.Range("N134").Value = -15208.2
$ws.Range("H136").Value = 1626
$ws.Range("I136").Value = 1145.8667
$ws.Range("K136").Value = 3437.6001
$ws.Range("M136").Value = -887.6001000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2601.43
$ws.Range("I4").Value = 2617.606
$ws.Range("K4").Value = 7852.818000000001
$ws.Range("M4").Value = -7740.818000000001
$ws.Range("H68").Value = 6671550
$ws.Range("I68").Value = 1037.2858
$ws.Range("J68").Value = 12508249
$ws.Range("K68").Value = 3111.8574
$ws.Range("L68").Value = 37524747
$ws.Range("M68").Value = -2300.8574
$ws.Range("N68").Value = -37526369
$ws.Range("H71").Value = 6671550
$ws.Range("I71").Value = 1037.2858
$ws.Range("J71").Value = 12508249
$ws.Range("K71").Value = 9335.572200000001
$ws.Range("L71").Value = 112574241
$ws.Range("M71").Value = -5279.572200000001
$ws.Range("N71").Value = -112582353
$ws.Range("H97").Value = 1003199.8
$ws.Range("J97").Value = 3999.75
$ws.Range("L97").Value = 11999.25
$ws.Range("N97").Value = -12991.25
$ws.Range("H117").Value = 1765.25
$ws.Range("J117").Value = 1618
$ws.Range("L117").Value = 4854
$ws.Range("N117").Value = -11738
$ws.Range("H121").Value = 9191442
$ws.Range("I121").Value = 20000492
$ws.Range("J121").Value = 183900.5
$ws.Range("K121").Value = 60001476
$ws.Range("L121").Value = 551701.5
$ws.Range("M121").Value = -60000166
$ws.Range("N121").Value = -554321.5
$ws.Range("H132").Value = 1916.8572
$ws.Range("I132").Value = 1278.25
$ws.Range("J132").Value = 2768.3333
$ws.Range("K132").Value = 11504.25
$ws.Range("L132").Value = 24914.9997
$ws.Range("M132").Value = -8974.25
$ws.Range("N132").Value = -29974.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5124.2144
$ws.Range("I102").Value = 2328
$ws.Range("J102").Value = 5459.76
$ws.Range("K102").Value = 2328
$ws.Range("L102").Value = 5459.76
$ws.Range("M102").Value = -706
$ws.Range("N102").Value = -8703.76

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1333.5
$ws.Range("I16").Value = 1333.5
$ws.Range("K16").Value = 1333.5
$ws.Range("M16").Value = -1163.5
$ws.Range("H55").Value = 439.27777
$ws.Range("I55").Value = 237
$ws.Range("J55").Value = 757.1429000000001
$ws.Range("K55").Value = 237
$ws.Range("L55").Value = 757.1429000000001
$ws.Range("M55").Value = -64
$ws.Range("N55").Value = -1103.1429
$ws.Range("H93").Value = 2871.2307
$ws.Range("I93").Value = 2714.4443
$ws.Range("K93").Value = 2714.4443
$ws.Range("M93").Value = -1466.4443
$ws.Range("H100").Value = 2947.25
$ws.Range("I100").Value = 2999.5
$ws.Range("J100").Value = 2895
$ws.Range("K100").Value = 2999.5
$ws.Range("L100").Value = 2895
$ws.Range("M100").Value = -2458.5
$ws.Range("N100").Value = -3977
$ws.Range("H136").Value = 5328.7744
$ws.Range("I136").Value = 4745.5
$ws.Range("K136").Value = 14236.5
$ws.Range("M136").Value = -11686.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3295.3215
$ws.Range("I132").Value = 3315.9583
$ws.Range("K132").Value = 9947.874899999999
$ws.Range("M132").Value = -7417.874899999999
